# Applies the "v1 improvements" commit to PlantillaAppsheetWoocommerce.xlsx:
#  - expands the "orders" sheet with the full WooCommerce order field set
#  - adds a new "order_refunds" sheet (becomes the active tab)
#  - leaves a selection of H17 on "orders" (as captured in the author's session)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "orders" sheet: rebuild row 1 headers (A1:V1), re-using the original
#    strings (id, customer_name, total, status, date_created) and inserting
#    all of the new WooCommerce order columns around them.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("orders")

$ordersHeaders = @(
    "id",
    "customer_name",
    "address_1",
    "address_2",
    "city",
    "postcode",
    "country",
    "email",
    "phone",
    "currency",
    "payment_method",
    "payment_method_title",
    "discount_total",
    "discount_tax",
    "shipping_total",
    "shipping_tax",
    "cart_tax",
    "total",
    "total_tax",
    "status",
    "date_created",
    "date_modified"
)

for ($i = 0; $i -lt $ordersHeaders.Length; $i++) {
    $ws1.Cells.Item(1, $i + 1).Value = $ordersHeaders[$i]
}

# Matches the saved session's lingering selection on the "orders" sheet.
[void]$ws1.Range("H17").Select()

# ---------------------------------------------------------------------------
# 2) Add the new "order_refunds" sheet after "order_details" and populate
#    its header row. Adding it last makes it the active/selected tab, which
#    mirrors activeTab="2" / tabSelected on sheet3 in the target workbook.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws3.Name = "order_refunds"

$refundsHeaders = @("id", "order_id", "reason", "total")
for ($i = 0; $i -lt $refundsHeaders.Length; $i++) {
    $ws3.Cells.Item(1, $i + 1).Value = $refundsHeaders[$i]
}
